$d = $word.ActiveDocument

# --- First paragraph: the "**ID__...__ID**" placeholder line ---
$p1 = $d.Paragraphs(1)
$pf = $p1.Format

# Add a paragraph border (top/left/bottom/right) with 5pt space, no line
# (mirrors <w:pBdr><w:top w:space="5"/>...</w:pBdr>).
$pf.Borders.DistanceFromTop = 5
$pf.Borders.DistanceFromLeft = 5
$pf.Borders.DistanceFromBottom = 5
$pf.Borders.DistanceFromRight = 5

# Bump the left indent from 120 twips (6pt) to 225 twips (11.25pt).
$pf.LeftIndent = 11.25

# Locate the placeholder ID text that lives in the paragraph's first run.
$idRange = $d.Content
$idRange.Find.Execute("**ID__AFFARS_pgi_5335_topic_7__ID**") | Out-Null

# The very next character is the standalone trailing-space run; drop it.
$spaceRange = $d.Range($idRange.End, $idRange.End + 1)
$spaceRange.Text = ""

# Update the ID text to the new topic id.
$idRange.Text = "**ID__AFFARS_SMC_PGI_5335__ID**"
